$wb = $excel.ActiveWorkbook

# --- Sheet 1: Fixture frances ---
$ws1 = $wb.Worksheets.Item("Fixture frances")
$ws1.Cells.Item(2, 2).Value = "@VEN"
$ws1.Cells.Item(2, 3).Value = "COL"
$ws1.Cells.Item(2, 4).Value = "ARG"
$ws1.Cells.Item(2, 5).Value = "@CHI"
$ws1.Cells.Item(2, 6).Value = "@PER"
$ws1.Cells.Item(2, 7).Value = "PAR"
$ws1.Cells.Item(2, 8).Value = "@URU"
$ws1.Cells.Item(2, 9).Value = "ECU"
$ws1.Cells.Item(2, 10).Value = "BOL"
$ws1.Cells.Item(2, 11).Value = "@COL"
$ws1.Cells.Item(2, 12).Value = "@ARG"
$ws1.Cells.Item(2, 13).Value = "CHI"
$ws1.Cells.Item(2, 14).Value = "PER"
$ws1.Cells.Item(2, 15).Value = "@PAR"
$ws1.Cells.Item(2, 16).Value = "URU"
$ws1.Cells.Item(2, 17).Value = "@ECU"
$ws1.Cells.Item(2, 18).Value = "@BOL"
$ws1.Cells.Item(2, 19).Value = "VEN"
$ws1.Cells.Item(3, 2).Value = "ECU"
$ws1.Cells.Item(3, 3).Value = "@URU"
$ws1.Cells.Item(3, 4).Value = "@BRA"
$ws1.Cells.Item(3, 5).Value = "COL"
$ws1.Cells.Item(3, 6).Value = "BOL"
$ws1.Cells.Item(3, 7).Value = "@VEN"
$ws1.Cells.Item(3, 8).Value = "CHI"
$ws1.Cells.Item(3, 9).Value = "@PAR"
$ws1.Cells.Item(3, 10).Value = "@PER"
$ws1.Cells.Item(3, 11).Value = "URU"
$ws1.Cells.Item(3, 12).Value = "BRA"
$ws1.Cells.Item(3, 13).Value = "@COL"
$ws1.Cells.Item(3, 14).Value = "@BOL"
$ws1.Cells.Item(3, 15).Value = "VEN"
$ws1.Cells.Item(3, 16).Value = "@CHI"
$ws1.Cells.Item(3, 17).Value = "PAR"
$ws1.Cells.Item(3, 18).Value = "PER"
$ws1.Cells.Item(3, 19).Value = "@ECU"
$ws1.Cells.Item(4, 2).Value = "PER"
$ws1.Cells.Item(4, 3).Value = "@BRA"
$ws1.Cells.Item(4, 4).Value = "PAR"
$ws1.Cells.Item(4, 5).Value = "@ARG"
$ws1.Cells.Item(4, 6).Value = "@ECU"
$ws1.Cells.Item(4, 7).Value = "URU"
$ws1.Cells.Item(4, 8).Value = "@VEN"
$ws1.Cells.Item(4, 9).Value = "BOL"
$ws1.Cells.Item(4, 10).Value = "@CHI"
$ws1.Cells.Item(4, 11).Value = "BRA"
$ws1.Cells.Item(4, 12).Value = "@PAR"
$ws1.Cells.Item(4, 13).Value = "ARG"
$ws1.Cells.Item(4, 14).Value = "ECU"
$ws1.Cells.Item(4, 15).Value = "@URU"
$ws1.Cells.Item(4, 16).Value = "VEN"
$ws1.Cells.Item(4, 17).Value = "@BOL"
$ws1.Cells.Item(4, 18).Value = "CHI"
$ws1.Cells.Item(4, 19).Value = "@PER"
$ws1.Cells.Item(5, 2).Value = "@PAR"
$ws1.Cells.Item(5, 3).Value = "ARG"
$ws1.Cells.Item(5, 4).Value = "@BOL"
$ws1.Cells.Item(5, 5).Value = "ECU"
$ws1.Cells.Item(5, 6).Value = "CHI"
$ws1.Cells.Item(5, 7).Value = "@COL"
$ws1.Cells.Item(5, 8).Value = "BRA"
$ws1.Cells.Item(5, 9).Value = "@PER"
$ws1.Cells.Item(5, 10).Value = "VEN"
$ws1.Cells.Item(5, 11).Value = "@ARG"
$ws1.Cells.Item(5, 12).Value = "BOL"
$ws1.Cells.Item(5, 13).Value = "@ECU"
$ws1.Cells.Item(5, 14).Value = "@CHI"
$ws1.Cells.Item(5, 15).Value = "COL"
$ws1.Cells.Item(5, 16).Value = "@BRA"
$ws1.Cells.Item(5, 17).Value = "PER"
$ws1.Cells.Item(5, 18).Value = "@VEN"
$ws1.Cells.Item(5, 19).Value = "PAR"
$ws1.Cells.Item(6, 2).Value = "@BOL"
$ws1.Cells.Item(6, 3).Value = "PAR"
$ws1.Cells.Item(6, 4).Value = "@ECU"
$ws1.Cells.Item(6, 5).Value = "BRA"
$ws1.Cells.Item(6, 6).Value = "@URU"
$ws1.Cells.Item(6, 7).Value = "PER"
$ws1.Cells.Item(6, 8).Value = "@ARG"
$ws1.Cells.Item(6, 9).Value = "VEN"
$ws1.Cells.Item(6, 10).Value = "COL"
$ws1.Cells.Item(6, 11).Value = "@PAR"
$ws1.Cells.Item(6, 12).Value = "ECU"
$ws1.Cells.Item(6, 13).Value = "@BRA"
$ws1.Cells.Item(6, 14).Value = "URU"
$ws1.Cells.Item(6, 15).Value = "@PER"
$ws1.Cells.Item(6, 16).Value = "ARG"
$ws1.Cells.Item(6, 17).Value = "@VEN"
$ws1.Cells.Item(6, 18).Value = "@COL"
$ws1.Cells.Item(6, 19).Value = "BOL"
$ws1.Cells.Item(7, 2).Value = "@COL"
$ws1.Cells.Item(7, 3).Value = "BOL"
$ws1.Cells.Item(7, 4).Value = "VEN"
$ws1.Cells.Item(7, 5).Value = "@PAR"
$ws1.Cells.Item(7, 6).Value = "BRA"
$ws1.Cells.Item(7, 7).Value = "@CHI"
$ws1.Cells.Item(7, 8).Value = "@ECU"
$ws1.Cells.Item(7, 9).Value = "URU"
$ws1.Cells.Item(7, 10).Value = "ARG"
$ws1.Cells.Item(7, 11).Value = "@BOL"
$ws1.Cells.Item(7, 12).Value = "@VEN"
$ws1.Cells.Item(7, 13).Value = "PAR"
$ws1.Cells.Item(7, 14).Value = "@BRA"
$ws1.Cells.Item(7, 15).Value = "CHI"
$ws1.Cells.Item(7, 16).Value = "ECU"
$ws1.Cells.Item(7, 17).Value = "@URU"
$ws1.Cells.Item(7, 18).Value = "@ARG"
$ws1.Cells.Item(7, 19).Value = "COL"
$ws1.Cells.Item(8, 2).Value = "BRA"
$ws1.Cells.Item(8, 4).Value = "@PER"
$ws1.Cells.Item(8, 5).Value = "BOL"
$ws1.Cells.Item(8, 6).Value = "@PAR"
$ws1.Cells.Item(8, 7).Value = "ARG"
$ws1.Cells.Item(8, 8).Value = "COL"
$ws1.Cells.Item(8, 9).Value = "@CHI"
$ws1.Cells.Item(8, 10).Value = "@URU"
$ws1.Cells.Item(8, 12).Value = "PER"
$ws1.Cells.Item(8, 13).Value = "@BOL"
$ws1.Cells.Item(8, 14).Value = "PAR"
$ws1.Cells.Item(8, 15).Value = "@ARG"
$ws1.Cells.Item(8, 16).Value = "@COL"
$ws1.Cells.Item(8, 17).Value = "CHI"
$ws1.Cells.Item(8, 18).Value = "URU"
$ws1.Cells.Item(8, 19).Value = "@BRA"
$ws1.Cells.Item(9, 2).Value = "CHI"
$ws1.Cells.Item(9, 3).Value = "@PER"
$ws1.Cells.Item(9, 4).Value = "URU"
$ws1.Cells.Item(9, 5).Value = "@VEN"
$ws1.Cells.Item(9, 6).Value = "@ARG"
$ws1.Cells.Item(9, 7).Value = "ECU"
$ws1.Cells.Item(9, 8).Value = "PAR"
$ws1.Cells.Item(9, 9).Value = "@COL"
$ws1.Cells.Item(9, 10).Value = "@BRA"
$ws1.Cells.Item(9, 11).Value = "PER"
$ws1.Cells.Item(9, 12).Value = "@URU"
$ws1.Cells.Item(9, 13).Value = "VEN"
$ws1.Cells.Item(9, 14).Value = "ARG"
$ws1.Cells.Item(9, 15).Value = "@ECU"
$ws1.Cells.Item(9, 16).Value = "@PAR"
$ws1.Cells.Item(9, 17).Value = "COL"
$ws1.Cells.Item(9, 18).Value = "BRA"
$ws1.Cells.Item(9, 19).Value = "@CHI"
$ws1.Cells.Item(10, 2).Value = "URU"
$ws1.Cells.Item(10, 3).Value = "@CHI"
$ws1.Cells.Item(10, 4).Value = "@COL"
$ws1.Cells.Item(10, 5).Value = "PER"
$ws1.Cells.Item(10, 6).Value = "VEN"
$ws1.Cells.Item(10, 7).Value = "@BRA"
$ws1.Cells.Item(10, 8).Value = "@BOL"
$ws1.Cells.Item(10, 9).Value = "ARG"
$ws1.Cells.Item(10, 10).Value = "@ECU"
$ws1.Cells.Item(10, 11).Value = "CHI"
$ws1.Cells.Item(10, 12).Value = "COL"
$ws1.Cells.Item(10, 13).Value = "@PER"
$ws1.Cells.Item(10, 14).Value = "@VEN"
$ws1.Cells.Item(10, 15).Value = "BRA"
$ws1.Cells.Item(10, 16).Value = "BOL"
$ws1.Cells.Item(10, 17).Value = "@ARG"
$ws1.Cells.Item(10, 18).Value = "ECU"
$ws1.Cells.Item(10, 19).Value = "@URU"
$ws1.Cells.Item(11, 2).Value = "@ARG"
$ws1.Cells.Item(11, 4).Value = "CHI"
$ws1.Cells.Item(11, 5).Value = "@URU"
$ws1.Cells.Item(11, 6).Value = "COL"
$ws1.Cells.Item(11, 7).Value = "@BOL"
$ws1.Cells.Item(11, 8).Value = "PER"
$ws1.Cells.Item(11, 9).Value = "@BRA"
$ws1.Cells.Item(11, 10).Value = "PAR"
$ws1.Cells.Item(11, 12).Value = "@CHI"
$ws1.Cells.Item(11, 13).Value = "URU"
$ws1.Cells.Item(11, 14).Value = "@COL"
$ws1.Cells.Item(11, 15).Value = "BOL"
$ws1.Cells.Item(11, 16).Value = "@PER"
$ws1.Cells.Item(11, 17).Value = "BRA"
$ws1.Cells.Item(11, 18).Value = "@PAR"
$ws1.Cells.Item(11, 19).Value = "ARG"

# --- Sheet 2: Breaks y secuencias ---
$ws2 = $wb.Worksheets.Item("Breaks y secuencias")
$ws2.Cells.Item(4, 3).Value = 5
$ws2.Cells.Item(4, 4).Value = 4
$ws2.Cells.Item(6, 3).Value = 4
$ws2.Cells.Item(6, 4).Value = 5
$ws2.Cells.Item(7, 3).Value = 4
$ws2.Cells.Item(7, 4).Value = 5
$ws2.Cells.Item(10, 3).Value = 5
$ws2.Cells.Item(10, 4).Value = 4

# --- Sheet 4: Partidos acumulados ---
$ws4 = $wb.Worksheets.Item("Partidos acumulados")
$ws4.Cells.Item(2, 6).Value = 3
$ws4.Cells.Item(2, 8).Value = 4
$ws4.Cells.Item(2, 14).Value = 6
$ws4.Cells.Item(2, 16).Value = 7
$ws4.Cells.Item(3, 4).Value = 2
$ws4.Cells.Item(3, 6).Value = 2
$ws4.Cells.Item(3, 8).Value = 3
$ws4.Cells.Item(3, 12).Value = 5
$ws4.Cells.Item(3, 14).Value = 7
$ws4.Cells.Item(3, 16).Value = 8
$ws4.Cells.Item(4, 2).Value = 0
$ws4.Cells.Item(4, 4).Value = 1
$ws4.Cells.Item(4, 6).Value = 3
$ws4.Cells.Item(4, 8).Value = 4
$ws4.Cells.Item(4, 10).Value = 5
$ws4.Cells.Item(4, 12).Value = 6
$ws4.Cells.Item(4, 14).Value = 6
$ws4.Cells.Item(4, 16).Value = 7
$ws4.Cells.Item(4, 18).Value = 8
$ws4.Cells.Item(5, 6).Value = 2
$ws4.Cells.Item(5, 8).Value = 3
$ws4.Cells.Item(5, 14).Value = 7
$ws4.Cells.Item(5, 16).Value = 8
$ws4.Cells.Item(6, 2).Value = 1
$ws4.Cells.Item(6, 10).Value = 4
$ws4.Cells.Item(6, 18).Value = 9
$ws4.Cells.Item(7, 2).Value = 1
$ws4.Cells.Item(7, 4).Value = 1
$ws4.Cells.Item(7, 8).Value = 4
$ws4.Cells.Item(7, 10).Value = 4
$ws4.Cells.Item(7, 12).Value = 6
$ws4.Cells.Item(7, 16).Value = 7
$ws4.Cells.Item(7, 18).Value = 9
$ws4.Cells.Item(8, 4).Value = 2
$ws4.Cells.Item(8, 12).Value = 5
$ws4.Cells.Item(9, 6).Value = 3
$ws4.Cells.Item(9, 8).Value = 3
$ws4.Cells.Item(9, 14).Value = 6
$ws4.Cells.Item(9, 16).Value = 8
$ws4.Cells.Item(10, 2).Value = 0
$ws4.Cells.Item(10, 10).Value = 5
$ws4.Cells.Item(10, 18).Value = 8
$ws4.Cells.Item(11, 6).Value = 2
$ws4.Cells.Item(11, 14).Value = 7
